# "Kosten-Nutzen Analyse" update:
#  - correct the "Übertragungsqualität" score for Cisco (G5) 5 -> 6
#    (dependent formulas / sums recalc automatically)
#  - highlight the three "Summe/100" result cells with a traffic-light
#    color scheme: best score green, worst score red, middle yellow
#  - move the active cell selection to G6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data correction -------------------------------------------------
$ws.Range("G5").Value = 6

# --- traffic-light colouring of the "Summe/100" totals ---------------
# D10 = 3CX total (highest -> green), H10 = Cisco total (middle -> yellow),
# F10 = Asterisk total (lowest -> red)
$ws.Range("D10").Interior.Color = 5296274   # RGB(146,208,80) - green
$ws.Range("H10").Interior.Color = 65535     # RGB(255,255,0)  - yellow
$ws.Range("F10").Interior.Color = 255       # RGB(255,0,0)    - red

# --- selection ---------------------------------------------------------
[void]$ws.Range("G6").Select()
